$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.9
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 2.5
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.4
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 13
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 7.5
$ws.Range("AH2").Value = 17
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 13
$ws.Range("AL2").Value = 26
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 10
$ws.Range("AP2").Value = 17
$ws.Range("AQ2").Value = 29
$ws.Range("AS2").Value = 101
$ws.Range("AT2").Value = 3.5
$ws.Range("AW2").Value = 351
$ws.Range("AX2").Value = 6
$ws.Range("AY2").Value = 19
$ws.Range("AZ2").Value = 21
$ws.Range("I3").Value = 2.8
$ws.Range("J3").Value = 2.82
$ws.Range("O3").Value = 1.24
$ws.Range("P3").Value = 3.62
$ws.Range("Q3").Value = 1.72
$ws.Range("R3").Value = 1.88
$ws.Range("S3").Value = 1.35
$ws.Range("T3").Value = 2.94
$ws.Range("U3").Value = 1.6
$ws.Range("V3").Value = 2.2
$ws.Range("W3").Value = 7.7
$ws.Range("X3").Value = 10.5
$ws.Range("Y3").Value = 7.6
$ws.Range("Z3").Value = 20
$ws.Range("AA3").Value = 14.5
$ws.Range("AB3").Value = 19
$ws.Range("AC3").Value = 10.5
$ws.Range("AE3").Value = 9.75
$ws.Range("AF3").Value = 35
$ws.Range("AK3").Value = 28
$ws.Range("AL3").Value = 18
$ws.Range("AN3").Value = 4.45
$ws.Range("AP3").Value = 17
$ws.Range("AQ3").Value = 45
$ws.Range("AR3").Value = 65
$ws.Range("AS3").Value = 175
$ws.Range("AT3").Value = 2.95
$ws.Range("AU3").Value = 6.3
$ws.Range("AV3").Value = 45
$ws.Range("AX3").Value = 5
$ws.Range("AY3").Value = 15
$ws.Range("BA3").Value = 65
$ws.Range("G4").Value = 1.52
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 5.4
$ws.Range("J4").Value = 2.05
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 5.4
$ws.Range("M4").Value = 1.02
$ws.Range("N4").Value = 9.699999999999999
$ws.Range("O4").Value = 1.24
$ws.Range("P4").Value = 3.62
$ws.Range("Q4").Value = 1.78
$ws.Range("R4").Value = 1.82
$ws.Range("S4").Value = 1.34
$ws.Range("T4").Value = 2.99
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.87
$ws.Range("W4").Value = 5.6
$ws.Range("X4").Value = 6
$ws.Range("Z4").Value = 9
$ws.Range("AB4").Value = 21
$ws.Range("AC4").Value = 10.25
$ws.Range("AE4").Value = 14
$ws.Range("AF4").Value = 60
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 12
$ws.Range("AI4").Value = 27
$ws.Range("AJ4").Value = 14
$ws.Range("AK4").Value = 80
$ws.Range("AL4").Value = 45
$ws.Range("AM4").Value = 40
$ws.Range("AN4").Value = 3.35
$ws.Range("AO4").Value = 7.2
$ws.Range("AP4").Value = 17
$ws.Range("AQ4").Value = 22
$ws.Range("AR4").Value = 55
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.82
$ws.Range("AU4").Value = 7.6
$ws.Range("AV4").Value = 70
$ws.Range("AX4").Value = 7.1
$ws.Range("AY4").Value = 32
$ws.Range("AZ4").Value = 35
$ws.Range("BA4").Value = 200
$ws.Range("BB4").Value = 200
$ws.Range("BC4").Value = 450
$ws.Range("H5").Value = 3.5
$ws.Range("L5").Value = 3.2
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 4.5
$ws.Range("Q5").Value = 1.62
$ws.Range("R5").Value = 2.25
$ws.Range("S5").Value = 1.3
$ws.Range("T5").Value = 3.4
$ws.Range("U5").Value = 1.53
$ws.Range("V5").Value = 2.38
$ws.Range("AB5").Value = 21
$ws.Range("AC5").Value = 15
$ws.Range("AF5").Value = 34
$ws.Range("AG5").Value = 101
$ws.Range("AH5").Value = 12
$ws.Range("AL5").Value = 19
$ws.Range("AM5").Value = 23
$ws.Range("AP5").Value = 19
$ws.Range("AS5").Value = 101
$ws.Range("AT5").Value = 3.4
$ws.Range("AU5").Value = 7
$ws.Range("BC5").Value = 101
$ws.Range("G8").Value = 1.55
$ws.Range("H8").Value = 3.9
$ws.Range("J8").Value = 2.02
$ws.Range("K8").Value = 2.3
$ws.Range("Q8").Value = 1.53
$ws.Range("R8").Value = 2.18
$ws.Range("S8").Value = 1.27
$ws.Range("T8").Value = 3.42
$ws.Range("U8").Value = 1.57
$ws.Range("V8").Value = 2.1
$ws.Range("W8").Value = 9
$ws.Range("X8").Value = 9
$ws.Range("Y8").Value = 7.9
$ws.Range("Z8").Value = 12.5
$ws.Range("AA8").Value = 11.25
$ws.Range("AB8").Value = 19
$ws.Range("AC8").Value = 14
$ws.Range("AD8").Value = 7.9
$ws.Range("AF8").Value = 50
$ws.Range("AG8").Value = 300
$ws.Range("AH8").Value = 17
$ws.Range("AN8").Value = 3.6
$ws.Range("AO8").Value = 7.2
$ws.Range("AP8").Value = 13.5
$ws.Range("AQ8").Value = 21
$ws.Range("AT8").Value = 3.2
$ws.Range("AV8").Value = 50
$ws.Range("AX8").Value = 7.1
$ws.Range("AZ8").Value = 29
$ws.Range("BB8").Value = 175
$ws.Range("BC8").Value = 300
$ws.Range("G9").Value = 1.95
$ws.Range("I9").Value = 4.5
$ws.Range("J9").Value = 2.75
$ws.Range("M9").Value = 1.1
$ws.Range("O9").Value = 1.44
$ws.Range("P9").Value = 2.63
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.5
$ws.Range("S9").Value = 1.57
$ws.Range("T9").Value = 2.25
$ws.Range("X9").Value = 8
$ws.Range("Z9").Value = 17
$ws.Range("AD9").Value = 6
$ws.Range("AT9").Value = 2.25
$ws.Range("G10").Value = 1.83
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 3.8
$ws.Range("J10").Value = 2.4
$ws.Range("L10").Value = 4.33
$ws.Range("M10").Value = 1.04
$ws.Range("O10").Value = 1.22
$ws.Range("U10").Value = 1.67
$ws.Range("V10").Value = 2.1
$ws.Range("Y10").Value = 8.5
$ws.Range("Z10").Value = 15
$ws.Range("AM10").Value = 34
$ws.Range("AO10").Value = 9.5
$ws.Range("AQ10").Value = 29
$ws.Range("AS10").Value = 101
$ws.Range("AX10").Value = 6
$ws.Range("AY10").Value = 21
$ws.Range("BA10").Value = 67
